$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E7").Value  = 16.547
$ws.Range("B8").Value  = 6.185
$ws.Range("B10").Value = 6.221
$ws.Range("B12").Value = 5.315
$ws.Range("E14").Value = 16.994
$ws.Range("E15").Value = 16.361
$ws.Range("B18").Value = 5.137
$ws.Range("E18").Value = 16.536
$ws.Range("E20").Value = 16.44600000000001
$ws.Range("B25").Value = 6.037
$ws.Range("E29").Value = 16.931
$ws.Range("E30").Value = 16.358
$ws.Range("E31").Value = 16.426
$ws.Range("E35").Value = 16.589
$ws.Range("B37").Value = 8.643000000000001
$ws.Range("E40").Value = 16.627
$ws.Range("E44").Value = 16.481
$ws.Range("E50").Value = 16.326
$ws.Range("E54").Value = 16.762
$ws.Range("B55").Value = 4.572
$ws.Range("B68").Value = 5.220000000000001
$ws.Range("E68").Value = 17.272
$ws.Range("E76").Value = 16.623
$ws.Range("B77").Value = 5.601
$ws.Range("B78").Value = 7.571000000000001
$ws.Range("B79").Value = 5.353
$ws.Range("B80").Value = 7.845999999999999
$ws.Range("B81").Value = 6.452
$ws.Range("B82").Value = 5.659000000000001
$ws.Range("B84").Value = 5.427000000000001
$ws.Range("E87").Value = 16.349
$ws.Range("E88").Value = 16.285
$ws.Range("E92").Value = 17.901
$ws.Range("E96").Value = 16.38
$ws.Range("E98").Value = 16.282
$ws.Range("B101").Value = 6.691
$ws.Range("E101").Value = 16.646
$ws.Range("B102").Value = 7.937
$ws.Range("E102").Value = 16.542
